# Update "Pagos" (column F) and "Inscrições homologadas" (column H) values
# for the rows that changed between the two data extracts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new "Pagos" (F) value
$newPagos = @{
    2  = 105
    3  = 44
    4  = 38
    5  = 140
    6  = 50
    7  = 42
    8  = 10
    10 = 631
    11 = 411
    12 = 667
    13 = 140
    14 = 116
    15 = 144
    16 = 179
    17 = 104
    18 = 47
    20 = 72
    21 = 114
    22 = 162
    23 = 173
    24 = 245
    25 = 284
    26 = 209
    27 = 313
    28 = 174
    29 = 147
    30 = 221
    31 = 53
    32 = 186
    33 = 245
    34 = 244
    35 = 167
    36 = 78
    37 = 160
    38 = 90
    39 = 148
    40 = 240
    41 = 350
    42 = 412
    43 = 120
    44 = 314
    45 = 156
    46 = 326
    47 = 472
    48 = 208
    49 = 239
    50 = 210
    51 = 174
    52 = 26
}

foreach ($row in $newPagos.Keys) {
    $pagos = $newPagos[$row]
    $isencoes = $ws.Cells.Item($row, 7).Value2
    $ws.Cells.Item($row, 6).Value2 = $pagos
    $ws.Cells.Item($row, 8).Value2 = $pagos + $isencoes
}
